# Update TPM-derived NATMI metrics for Wnt4-Fzd6 (recalculated with new TPM values)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.433753
$ws.Range("H2").Value = 4.301259
$ws.Range("I2").Value = 0.2048034833446371
$ws.Range("J2").Value = 0.2048034833446371
$ws.Range("M2").Value = 22.618885
$ws.Range("N2").Value = 67.856655
$ws.Range("O2").Value = 0.9027998993061069
$ws.Range("P2").Value = 0.902799899306107
$ws.Range("Q2").Value = 32.42989422540501
$ws.Range("R2").Value = 291.869048028645
$ws.Range("S2").Value = 0.1848965641410784
$ws.Range("T2").Value = 0.1848965641410784

# Row 3
$ws.Range("G3").Value = 1.433753
$ws.Range("H3").Value = 4.301259
$ws.Range("I3").Value = 0.2048034833446371
$ws.Range("J3").Value = 0.2048034833446371
$ws.Range("O3").Value = 0.08600330007856447
$ws.Range("P3").Value = 0.08600330007856449
$ws.Range("Q3").Value = 3.089364461302333
$ws.Range("R3").Value = 27.804280151721
$ws.Range("S3").Value = 0.01761377543522411
$ws.Range("T3").Value = 0.01761377543522411

# Row 4
$ws.Range("G4").Value = 1.433753
$ws.Range("H4").Value = 4.301259
$ws.Range("I4").Value = 0.2048034833446371
$ws.Range("J4").Value = 0.2048034833446371
$ws.Range("M4").Value = 0.2805263333333333
$ws.Range("N4").Value = 0.841579
$ws.Range("O4").Value = 0.0111968006153285
$ws.Range("P4").Value = 0.01119680061532851
$ws.Range("Q4").Value = 0.4022054719956666
$ws.Range("R4").Value = 3.619849247961
$ws.Range("S4").Value = 0.002293143768334654
$ws.Range("T4").Value = 0.002293143768334654

# Row 5
$ws.Range("H5").Value = 4.309093
$ws.Range("I5").Value = 0.2051764974989863
$ws.Range("J5").Value = 0.2051764974989863
$ws.Range("M5").Value = 22.618885
$ws.Range("N5").Value = 67.856655
$ws.Range("O5").Value = 0.9027998993061069
$ws.Range("P5").Value = 0.902799899306107
$ws.Range("Q5").Value = 32.48895967376833
$ws.Range("R5").Value = 292.400637063915
$ws.Range("S5").Value = 0.1852333212820645
$ws.Range("T5").Value = 0.1852333212820646

# Row 6
$ws.Range("H6").Value = 4.309093
$ws.Range("I6").Value = 0.2051764974989863
$ws.Range("J6").Value = 0.2051764974989863
$ws.Range("O6").Value = 0.08600330007856447
$ws.Range("P6").Value = 0.08600330007856449
$ws.Range("Q6").Value = 3.094991204818555
$ws.Range("R6").Value = 27.854920843367
$ws.Range("S6").Value = 0.01764585588347415
$ws.Range("T6").Value = 0.01764585588347416

# Row 7
$ws.Range("H7").Value = 4.309093
$ws.Range("I7").Value = 0.2051764974989863
$ws.Range("J7").Value = 0.2051764974989863
$ws.Range("M7").Value = 0.2805263333333333
$ws.Range("N7").Value = 0.841579
$ws.Range("O7").Value = 0.0111968006153285
$ws.Range("P7").Value = 0.01119680061532851
$ws.Range("Q7").Value = 0.4029380197607778
$ws.Range("R7").Value = 3.626442177847
$ws.Range("S7").Value = 0.002297320333447597
$ws.Range("T7").Value = 0.002297320333447598

# Row 8
$ws.Range("G8").Value = 2.046430333333333
$ws.Range("H8").Value = 6.139291
$ws.Range("I8").Value = 0.2923209651096064
$ws.Range("J8").Value = 0.2923209651096064
$ws.Range("M8").Value = 22.618885
$ws.Range("N8").Value = 67.856655
$ws.Range("O8").Value = 0.9027998993061069
$ws.Range("P8").Value = 0.902799899306107
$ws.Range("Q8").Value = 46.28797237017834
$ws.Range("R8").Value = 416.5917513316051
$ws.Range("S8").Value = 0.2639073378660167
$ws.Range("T8").Value = 0.2639073378660167

# Row 9
$ws.Range("G9").Value = 2.046430333333333
$ws.Range("H9").Value = 6.139291
$ws.Range("I9").Value = 0.2923209651096064
$ws.Range("J9").Value = 0.2923209651096064
$ws.Range("O9").Value = 0.08600330007856447
$ws.Range("P9").Value = 0.08600330007856449
$ws.Range("Q9").Value = 4.409524614303222
$ws.Range("R9").Value = 39.685721528729
$ws.Range("S9").Value = 0.02514056768157705
$ws.Range("T9").Value = 0.02514056768157706

# Row 10
$ws.Range("G10").Value = 2.046430333333333
$ws.Range("H10").Value = 6.139291
$ws.Range("I10").Value = 0.2923209651096064
$ws.Range("J10").Value = 0.2923209651096064
$ws.Range("M10").Value = 0.2805263333333333
$ws.Range("N10").Value = 0.841579
$ws.Range("O10").Value = 0.0111968006153285
$ws.Range("P10").Value = 0.01119680061532851
$ws.Range("Q10").Value = 0.5740775978321111
$ws.Range("R10").Value = 5.166698380489
$ws.Range("S10").Value = 0.003273059562012663
$ws.Range("T10").Value = 0.003273059562012664

# Row 11
$ws.Range("G11").Value = 0.6416706666666666
$ws.Range("H11").Value = 1.925012
$ws.Range("I11").Value = 0.0916590149721806
$ws.Range("J11").Value = 0.0916590149721806
$ws.Range("M11").Value = 22.618885
$ws.Range("N11").Value = 67.856655
$ws.Range("O11").Value = 0.9027998993061069
$ws.Range("P11").Value = 0.902799899306107
$ws.Range("Q11").Value = 14.51387501720667
$ws.Range("R11").Value = 130.62487515486
$ws.Range("S11").Value = 0.08274974948738159
$ws.Range("T11").Value = 0.0827497494873816

# Row 12
$ws.Range("G12").Value = 0.6416706666666666
$ws.Range("H12").Value = 1.925012
$ws.Range("I12").Value = 0.0916590149721806
$ws.Range("J12").Value = 0.0916590149721806
$ws.Range("O12").Value = 0.08600330007856447
$ws.Range("P12").Value = 0.08600330007856449
$ws.Range("Q12").Value = 1.382633238403111
$ws.Range("R12").Value = 12.443699145628
$ws.Range("S12").Value = 0.007882977769558082
$ws.Range("T12").Value = 0.007882977769558083

# Row 13
$ws.Range("G13").Value = 0.6416706666666666
$ws.Range("H13").Value = 1.925012
$ws.Range("I13").Value = 0.0916590149721806
$ws.Range("J13").Value = 0.0916590149721806
$ws.Range("M13").Value = 0.2805263333333333
$ws.Range("N13").Value = 0.841579
$ws.Range("O13").Value = 0.0111968006153285
$ws.Range("P13").Value = 0.01119680061532851
$ws.Range("Q13").Value = 0.1800055193275555
$ws.Range("R13").Value = 1.620049673948
$ws.Range("S13").Value = 0.001026287715240916
$ws.Range("T13").Value = 0.001026287715240916

# Row 14
$ws.Range("G14").Value = 1.442409666666667
$ws.Range("H14").Value = 4.327229
$ws.Range("I14").Value = 0.2060400390745897
$ws.Range("J14").Value = 0.2060400390745897
$ws.Range("M14").Value = 22.618885
$ws.Range("N14").Value = 67.856655
$ws.Range("O14").Value = 0.9027998993061069
$ws.Range("P14").Value = 0.902799899306107
$ws.Range("Q14").Value = 32.62569837322167
$ws.Range("R14").Value = 293.631285358995
$ws.Range("S14").Value = 0.1860129265295659
$ws.Range("T14").Value = 0.186012926529566

# Row 15
$ws.Range("G15").Value = 1.442409666666667
$ws.Range("H15").Value = 4.327229
$ws.Range("I15").Value = 0.2060400390745897
$ws.Range("J15").Value = 0.2060400390745897
$ws.Range("O15").Value = 0.08600330007856447
$ws.Range("P15").Value = 0.08600330007856449
$ws.Range("Q15").Value = 3.108017324350111
$ws.Range("R15").Value = 27.972155919151
$ws.Range("S15").Value = 0.01772012330873109
$ws.Range("T15").Value = 0.01772012330873109

# Row 16
$ws.Range("G16").Value = 1.442409666666667
$ws.Range("H16").Value = 4.327229
$ws.Range("I16").Value = 0.2060400390745897
$ws.Range("J16").Value = 0.2060400390745897
$ws.Range("M16").Value = 0.2805263333333333
$ws.Range("N16").Value = 0.841579
$ws.Range("O16").Value = 0.0111968006153285
$ws.Range("P16").Value = 0.01119680061532851
$ws.Range("Q16").Value = 0.4046338949545555
$ws.Range("R16").Value = 3.641705054591
$ws.Range("S16").Value = 0.002306989236292675
$ws.Range("T16").Value = 0.002306989236292676
